# Re-pull / push updated dSF (column F) values following a repull of source data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = 1
    7  = 4
    8  = 1
    11 = -2
    12 = 9
    19 = 3
    20 = -3
    28 = 3
    35 = -1
    39 = -1
    42 = -2
    43 = -2
    50 = -4
    53 = 6
    56 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
